# Fruta / hortaliza, semanal
#
# A new weekly price record for "Haba" (Mercado Mayorista Lo Valledor de
# Santiago) needs to be inserted right before the existing row that used to
# be row 405 (dated 2021-12-02 / serial 44532). Inserting a whole row shifts
# every following record down by one (old row 405 -> new row 406, ...,
# old row 454 -> new row 455), which matches the diff exactly and bumps the
# sheet's used range from A1:R454 to A1:R455.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 405..454 down to 406..455, opening up a blank row 405.
$ws.Rows.Item(405).Insert()

# Populate the newly inserted row 405 with the new record. Columns A, B, C,
# E, F, G, H, I, N, Q, R repeat the same market/category metadata as the
# neighbouring rows; D, J, K, L, M, O, P carry the new observation's values.
$ws.Range("A405").Value = 6
$ws.Range("B405").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C405").Value = 'Metropolitana'
$ws.Range("D405").Value = 45194
$ws.Range("E405").Value = 13
$ws.Range("F405").Value = 100112026
$ws.Range("G405").Value = 'Haba'
$ws.Range("H405").Value = 'Sin especificar'
$ws.Range("I405").Value = 'Primera'
$ws.Range("J405").Value = 350
$ws.Range("K405").Value = 10000
$ws.Range("L405").Value = 12000
$ws.Range("M405").Value = 11314
$ws.Range("N405").Value = '$/saco 25 kilos'
$ws.Range("O405").Value = 'Provincia de Melipilla'
$ws.Range("P405").Value = 453
$ws.Range("Q405").Value = 25
$ws.Range("R405").Value = 'Hortaliza'
